$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a Price-column (D) value while forcing the cell to be
# stored as plain text. This keeps numeric-looking strings (e.g.
# "298.99", "0.120") with their exact literal representation (trailing
# zeros, thousand-dot groupings, etc.) instead of letting Excel's normal
# type inference silently convert them to floating point numbers.
function Set-PriceText($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "42.581.33"
$ws.Range("E2").Value = "  -7.47%  "

# Row 3 - Ethereum
Set-PriceText "D3" "2.544.41"
$ws.Range("E3").Value = "  -3.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-PriceText "D5" "298.99"
$ws.Range("E5").Value = "  -3.75%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  -6.58%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -4.29%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
Set-PriceText "D9" "0.549"
$ws.Range("E9").Value = "  -5.63%  "

# Row 10 - Avalanche
Set-PriceText "D10" "35.86"
$ws.Range("E10").Value = "  -7.30%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -5.17%  "

# Row 12 - Polkadot
Set-PriceText "D12" "7.66"
$ws.Range("E12").Value = "  -5.36%  "

# Row 13 - TRON
Set-PriceText "D13" "0.114"
$ws.Range("E13").Value = "  +5.62%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-PriceText "D14" "2.934.14"

# Row 15 - WrappedEther
Set-PriceText "D15" "2.542.81"
$ws.Range("E15").Value = "  -3.67%  "

# Row 16 - Polygon
Set-PriceText "D16" "0.873"
$ws.Range("E16").Value = "  -5.65%  "

# Row 17 - Chainlink
$ws.Range("E17").Value = "  -4.86%  "

# Row 18 - WrappedBTC
Set-PriceText "D18" "42.613.71"
$ws.Range("E18").Value = "  -7.48%  "

# Row 19
Set-PriceText "D19" "12.85"
$ws.Range("E19").Value = "  -0.28%  "

# Row 20
$ws.Range("E20").Value = "  -3.82%  "

# Row 21
Set-PriceText "D21" "6.54"
$ws.Range("E21").Value = "  -3.62%  "

# Row 22
Set-PriceText "D22" "71.37"
$ws.Range("E22").Value = "  -4.29%  "

# Row 23
Set-PriceText "D23" "255.62"
$ws.Range("E23").Value = "  -9.46%  "

# Row 24
$ws.Range("E24").Value = "  -5.19%  "

# Row 25
$ws.Range("E25").Value = "  -4.87%  "

# Row 26
$ws.Range("E26").Value = "  -6.04%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.01%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -5.13%  "

# Row 29 - InjectiveProtocol
Set-PriceText "D29" "36.88"
$ws.Range("E29").Value = "  -4.51%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -5.55%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  -5.18%  "

# Row 32 - Monero
Set-PriceText "D32" "152.73"
$ws.Range("E32").Value = "  -2.29%  "

# Row 33 - ARBITRUM
Set-PriceText "D33" "2.16"
$ws.Range("E33").Value = "  -7.72%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  -2.29%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -8.79%  "

# Row 36 - Hedera
Set-PriceText "D36" "0.0792"
$ws.Range("E36").Value = "  -5.94%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -6.88%  "

# Row 38 - was Stellar, now Celestia
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-PriceText "D38" "17.15"
$ws.Range("E38").Value = "  +8.46%  "

# Row 39 - was Celestia, now Stellar
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-PriceText "D39" "0.120"
$ws.Range("E39").Value = "  -3.71%  "

# Row 40 - EnergySwap
Set-PriceText "D40" "23.91"
$ws.Range("E40").Value = "  -8.60%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  -5.91%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  -4.23%  "

# Row 43 - NEARProtocol
$ws.Range("E43").Value = "  -5.40%  "

# Row 44 - Maker
Set-PriceText "D44" "2.087.12"
$ws.Range("E44").Value = "  -3.29%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  -0.06%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -1.24%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +2.72%  "

# Row 48 - BitcoinSV
Set-PriceText "D48" "84.33"
$ws.Range("E48").Value = "  -10.25%  "

# Row 49 - RocketPoolETH
Set-PriceText "D49" "2.789.94"
$ws.Range("E49").Value = "  -4.00%  "

# Row 50 - Aave
Set-PriceText "D50" "103.95"
$ws.Range("E50").Value = "  -6.61%  "

# Row 51 - Stacks
$ws.Range("E51").Value = "  -5.28%  "
